$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1650
$ws.Range("F3").Value = 863
$ws.Range("F4").Value = 273
$ws.Range("F5").Value = 83
$ws.Range("F6").Value = 1186
$ws.Range("F7").Value = 800
$ws.Range("F8").Value = 827
$ws.Range("F9").Value = 1529
$ws.Range("F10").Value = 306
$ws.Range("F11").Value = 1061
$ws.Range("F12").Value = 32
$ws.Range("F14").Value = 202
$ws.Range("F15").Value = 60
$ws.Range("F16").Value = 511
$ws.Range("F17").Value = 64
$ws.Range("F18").Value = 41
$ws.Range("F22").Value = 580
$ws.Range("F23").Value = 584
$ws.Range("F24").Value = 54
$ws.Range("F26").Value = 781
$ws.Range("F27").Value = 260
$ws.Range("F28").Value = 195
$ws.Range("F29").Value = 1
$ws.Range("F30").Value = 375

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 11
$ws.Range("F3").Value = 1035
$ws.Range("F5").Value = 281
$ws.Range("F6").Value = 18
$ws.Range("F7").Value = 152
$ws.Range("F8").Value = 70
$ws.Range("F10").Value = 89
$ws.Range("F11").Value = 14

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 266

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 266
$ws.Range("F3").Value = 1650
$ws.Range("F4").Value = 11
$ws.Range("F5").Value = 863
$ws.Range("F6").Value = 273
$ws.Range("F7").Value = 1035
$ws.Range("F8").Value = 83
$ws.Range("F9").Value = 1186
$ws.Range("F10").Value = 800
$ws.Range("F11").Value = 827
$ws.Range("F12").Value = 1529
$ws.Range("F13").Value = 306
$ws.Range("F14").Value = 1061
$ws.Range("F15").Value = 32
$ws.Range("F17").Value = 202
$ws.Range("F18").Value = 60
$ws.Range("F19").Value = 511
$ws.Range("F20").Value = 64
$ws.Range("F21").Value = 41
$ws.Range("F24").Value = 281
$ws.Range("F27").Value = 18
$ws.Range("F28").Value = 152
$ws.Range("F29").Value = 152
$ws.Range("F30").Value = 580
$ws.Range("F31").Value = 584
$ws.Range("F32").Value = 54
$ws.Range("F34").Value = 781
$ws.Range("F35").Value = 260
$ws.Range("F36").Value = 70
$ws.Range("F37").Value = 195
$ws.Range("F39").Value = 89
$ws.Range("F40").Value = 89
$ws.Range("F41").Value = 1
$ws.Range("F42").Value = 14
$ws.Range("F43").Value = 375

